$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.970.52'
$ws.Range('E2').Value = '  +3.09%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.909.58'
$ws.Range('E3').Value = '  +1.43%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.31%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '245.21'
$ws.Range('E5').Value = '  +0.43%  '

$ws.Range('E6').Value = '  +0.26%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5007'
$ws.Range('E7').Value = '  +0.98%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2992'
$ws.Range('E8').Value = '  +2.45%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06867'

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.910.47'
$ws.Range('E10').Value = '  +1.50%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '16.99'
$ws.Range('E11').Value = '  +0.33%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07303'
$ws.Range('E12').Value = '  +1.45%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '91.32'
$ws.Range('E13').Value = '  +6.21%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.096'
$ws.Range('E14').Value = '  +5.22%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6797'
$ws.Range('E15').Value = '  +1.88%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '30.949.14'
$ws.Range('E16').Value = '  +3.10%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000008031'
$ws.Range('E17').Value = '  +2.76%  '

$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  +0.35%  '

$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.28'
$ws.Range('E19').Value = '  +3.58%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '2.159.45'
$ws.Range('E20').Value = '  +1.81%  '

$ws.Range('E21').Value = '  +0.28%  '

$ws.Range('E22').Value = '  +2.37%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '183.69'
$ws.Range('E23').Value = '  +35.27%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.117'
$ws.Range('E24').Value = '  +9.21%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.367'
$ws.Range('E25').Value = '  +2.25%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '154.15'
$ws.Range('E26').Value = '  +2.69%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.92'
$ws.Range('E27').Value = '  +12.95%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.940'
$ws.Range('E28').Value = '  +1.81%  '

$ws.Range('E30').Value = '  +4.09%  '

$ws.Range('E31').Value = '  +3.45%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.053'
$ws.Range('E32').Value = '  +2.60%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05250'
$ws.Range('E33').Value = '  +5.26%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7505'
$ws.Range('E34').Value = '  +6.84%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.140'
$ws.Range('E35').Value = '  +3.19%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.670'
$ws.Range('E36').Value = '  +0.67%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.01945'
$ws.Range('E37').Value = '  +18.50%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.736'
$ws.Range('E38').Value = '  +1.54%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.183'
$ws.Range('E39').Value = '  -0.76%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9363'
$ws.Range('E40').Value = '  +0.11%  '

$ws.Range('E41').Value = '  +4.63%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '106.30'
$ws.Range('E42').Value = '  +4.72%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.879'
$ws.Range('E43').Value = '  -1.44%  '

$ws.Range('E44').Value = '  +0.13%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '7.810'
$ws.Range('E45').Value = '  +3.59%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.1345'
$ws.Range('E46').Value = '  +6.62%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.05846'
$ws.Range('E47').Value = '  +2.16%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.609'
$ws.Range('E48').Value = '  +4.85%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.3907'
$ws.Range('E49').Value = '  +5.43%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '33.27'
$ws.Range('E50').Value = '  +2.80%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.393'
$ws.Range('E51').Value = '  +4.26%  '
